$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Ventas objetivo" formula change -> recalculated values for
# "uds. Objetivo semana pasada" (R), "Diferencia Stock" (L),
# "Tendencia Consumo" (T), "Pedido Final" (U) and the summary metrics.

$ws.Range("R5").Value = 2

$ws.Range("L6").Value = 0
$ws.Range("R6").Value = 1

$ws.Range("R10").Value = 1
$ws.Range("T10").Value = 1

$ws.Range("R11").Value = 3

$ws.Range("R12").Value = 2

$ws.Range("R15").Value = 1
$ws.Range("T15").Value = 0

$ws.Range("R16").Value = 1

$ws.Range("R17").Value = 3

$ws.Range("L23").Value = 0
$ws.Range("R23").Value = 4

$ws.Range("R26").Value = 1

$ws.Range("R27").Value = 2
$ws.Range("T27").Value = 0

$ws.Range("R34").Value = 4

$ws.Range("R36").Value = 1
$ws.Range("T36").Value = 0

$ws.Range("R37").Value = 3
$ws.Range("T37").Value = 0

$ws.Range("L38").Value = 0
$ws.Range("R38").Value = 7
$ws.Range("T38").Value = 0
$ws.Range("U38").Value = 7

$ws.Range("L39").Value = 0
$ws.Range("R39").Value = 4

$ws.Range("C43").Value = 54

$ws.Range("C54").Value = 0
